$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H54").Value = 25499.5
$ws.Range("I54").Value = 25499.5
$ws.Range("K54").Value = 25499.5
$ws.Range("M54").Value = -25013.5

$ws.Range("H82").Value = 1492.6
$ws.Range("I82").Value = 1492.25
$ws.Range("K82").Value = 4476.75
$ws.Range("M82").Value = -4070.75

$ws.Range("H85").Value = 1492.6
$ws.Range("I85").Value = 1492.25
$ws.Range("K85").Value = 4476.75
$ws.Range("M85").Value = -3072.75

$ws.Range("H113").Value = 2899.4
$ws.Range("I113").Value = 2374.5
$ws.Range("K113").Value = 2374.5
$ws.Range("M113").Value = 879.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 5798.8945
$ws.Range("I45").Value = 3107.4546
$ws.Range("J45").Value = 9499.625
$ws.Range("K45").Value = 3107.4546
$ws.Range("L45").Value = 9499.625
$ws.Range("M45").Value = -2730.4546
$ws.Range("N45").Value = -10253.625

$ws.Range("H74").Value = 1056.22
$ws.Range("I74").Value = 1042.0698
$ws.Range("J74").Value = 1143.1428
$ws.Range("K74").Value = 1042.0698
$ws.Range("L74").Value = 1143.1428
$ws.Range("M74").Value = -168.0698
$ws.Range("N74").Value = -2891.1428

$ws.Range("H77").Value = 1056.22
$ws.Range("I77").Value = 1042.0698
$ws.Range("J77").Value = 1143.1428
$ws.Range("K77").Value = 5210.349
$ws.Range("L77").Value = 5715.714
$ws.Range("M77").Value = -842.3490000000002
$ws.Range("N77").Value = -14451.714

$ws.Range("H92").Value = 25025000
$ws.Range("J92").Value = 50000
$ws.Range("L92").Value = 50000
$ws.Range("N92").Value = -54992

$ws.Range("H97").Value = 1293.0769
$ws.Range("I97").Value = 1233.6666
$ws.Range("K97").Value = 1233.6666
$ws.Range("M97").Value = -737.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2632.6667
$ws.Range("I105").Value = 1449.5
$ws.Range("K105").Value = 1449.5
$ws.Range("M105").Value = 297.5

$ws.Range("H134").Value = 4029.6667
$ws.Range("I134").Value = 2583.5588
$ws.Range("J134").Value = 7541.643
$ws.Range("K134").Value = 7750.676399999999
$ws.Range("L134").Value = 22624.929
$ws.Range("M134").Value = -5215.676399999999
$ws.Range("N134").Value = -27694.929

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 479812.53
$ws.Range("I31").Value = 911027.6
$ws.Range("J31").Value = 5475.9
$ws.Range("K31").Value = 911027.6
$ws.Range("L31").Value = 5475.9
$ws.Range("M31").Value = -910732.6
$ws.Range("N31").Value = -6065.9

$ws.Range("H34").Value = 479812.53
$ws.Range("I34").Value = 911027.6
$ws.Range("J34").Value = 5475.9
$ws.Range("K34").Value = 911027.6
$ws.Range("L34").Value = 5475.9
$ws.Range("M34").Value = -910825.6
$ws.Range("N34").Value = -5879.9

$ws.Range("H58").Value = 241425.31
$ws.Range("I58").Value = 436290
$ws.Range("K58").Value = 436290
$ws.Range("M58").Value = -436087

$ws.Range("H62").Value = 4500
$ws.Range("I62").Value = 4500
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 4500
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -3876
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 4500
$ws.Range("I65").Value = 4500
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 22500
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -19380
$ws.Range("N65").ClearContents()

$ws.Range("H99").Value = 5041.857
$ws.Range("I99").Value = 3783.1428
$ws.Range("K99").Value = 3783.1428
$ws.Range("M99").Value = -2285.1428

$ws.Range("H105").Value = 1785.2354
$ws.Range("J105").Value = 1797.5
$ws.Range("L105").Value = 1797.5
$ws.Range("N105").Value = -5291.5

$ws.Range("H122").Value = 2974.4707
$ws.Range("I122").Value = 2154.3333
$ws.Range("J122").Value = 4942.8
$ws.Range("K122").Value = 6462.999899999999
$ws.Range("L122").Value = 14828.4
$ws.Range("M122").Value = -4012.999899999999
$ws.Range("N122").Value = -19728.4

$ws.Range("H126").Value = 5041.857
$ws.Range("I126").Value = 3783.1428
$ws.Range("K126").Value = 11349.4284
$ws.Range("M126").Value = -8879.428400000001

$ws.Range("H132").Value = 3415.1667
$ws.Range("I132").Value = 2426.6
$ws.Range("K132").Value = 7279.799999999999
$ws.Range("M132").Value = -4749.799999999999

$ws.Range("H134").Value = 3348.3877
$ws.Range("I134").Value = 2310.6562
$ws.Range("J134").Value = 5301.7646
$ws.Range("K134").Value = 6931.9686
$ws.Range("L134").Value = 15905.2938
$ws.Range("M134").Value = -4396.9686
$ws.Range("N134").Value = -20975.2938

$ws.Range("H136").Value = 241425.31
$ws.Range("I136").Value = 436290
$ws.Range("K136").Value = 1308870
$ws.Range("M136").Value = -1306320

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 3028.1667
$ws.Range("I132").Value = 433
$ws.Range("K132").Value = 3897
$ws.Range("M132").Value = -1367

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 23405
$ws.Range("J47").Value = 23405
$ws.Range("L47").Value = 23405
$ws.Range("N47").Value = -24541

$ws.Range("H132").Value = 775073.25
$ws.Range("I132").Value = 1431715.8
$ws.Range("J132").Value = 8990.333000000001
$ws.Range("K132").Value = 4295147.4
$ws.Range("L132").Value = 26970.999
$ws.Range("M132").Value = -4292617.4
$ws.Range("N132").Value = -32030.999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 5001.3335
$ws.Range("I20").Value = 5002
$ws.Range("K20").Value = 5002
$ws.Range("M20").Value = -4776

$ws.Range("H23").Value = 349998.44
$ws.Range("I23").Value = 504497.25
$ws.Range("K23").Value = 504497.25
$ws.Range("M23").Value = -504267.25

$ws.Range("H40").Value = 6192.5713
$ws.Range("I40").Value = 2939.4
$ws.Range("K40").Value = 2939.4
$ws.Range("M40").Value = -2803.4

$ws.Range("H55").Value = 1238.1111
$ws.Range("I55").Value = 522.8889
$ws.Range("J55").Value = 1953.3334
$ws.Range("K55").Value = 522.8889
$ws.Range("L55").Value = 1953.3334
$ws.Range("M55").Value = -349.8889
$ws.Range("N55").Value = -2299.3334

$ws.Range("H136").Value = 6381.5884
$ws.Range("I136").Value = 4779.6665
$ws.Range("J136").Value = 6724.857
$ws.Range("K136").Value = 14338.9995
$ws.Range("L136").Value = 20174.571
$ws.Range("M136").Value = -11788.9995
$ws.Range("N136").Value = -25274.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 99000
$ws.Range("J57").Value = 99000
$ws.Range("L57").Value = 99000
$ws.Range("N57").Value = -100508

$ws.Range("H92").Value = 107683.336
$ws.Range("J92").Value = 107683.336
$ws.Range("L92").Value = 107683.336
$ws.Range("N92").Value = -112675.336

$ws.Range("H126").Value = 3781.1765
$ws.Range("I126").Value = 2548.1
$ws.Range("K126").Value = 7644.299999999999
$ws.Range("M126").Value = -5174.299999999999

$ws.Range("H132").Value = 4448.625
$ws.Range("I132").Value = 1942.8572
$ws.Range("K132").Value = 5828.571599999999
$ws.Range("M132").Value = -3298.571599999999

$ws.Range("H136").Value = 1411.6666
$ws.Range("J136").Value = 2793.1428
$ws.Range("L136").Value = 8379.428400000001
$ws.Range("N136").Value = -13479.4284
